$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.727.37"
$ws.Range("E2").Value = "  +2.21%  "
$ws.Range("D3").Value = "1.892.44"
$ws.Range("E3").Value = "  +0.80%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.22"
$ws.Range("E5").Value = "  +0.56%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4924"
$ws.Range("E7").Value = "  +0.28%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2963"
$ws.Range("E8").Value = "  +1.11%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06792"
$ws.Range("E9").Value = "  +2.70%  "
$ws.Range("D10").Value = "1.888.86"
$ws.Range("E10").Value = "  +0.56%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "17.21"
$ws.Range("E11").Value = "  +3.44%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07241"
$ws.Range("E12").Value = "  +0.62%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "91.04"
$ws.Range("E13").Value = "  +5.36%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6787"
$ws.Range("E14").Value = "  +1.64%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.044"
$ws.Range("E15").Value = "  +2.73%  "
$ws.Range("D16").Value = "30.705.04"
$ws.Range("E16").Value = "  +2.27%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000007992"
$ws.Range("E17").Value = "  +2.41%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9998"
$ws.Range("E18").Value = "  +0.02%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.16"
$ws.Range("E19").Value = "  +2.89%  "
$ws.Range("D20").Value = "2.132.46"
$ws.Range("E20").Value = "  +0.49%  "
$ws.Range("E21").Value = "  +0.36%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.820"
$ws.Range("E22").Value = "  +0.69%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "188.59"
$ws.Range("E23").Value = "  +31.52%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.159"
$ws.Range("E24").Value = "  +4.76%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.372"
$ws.Range("E25").Value = "  +2.64%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "155.85"
$ws.Range("E26").Value = "  +2.23%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.06"
$ws.Range("E27").Value = "  +12.30%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.904"
$ws.Range("E28").Value = "  +0.25%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.401"
$ws.Range("E29").Value = "  +1.22%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.343"
$ws.Range("E30").Value = "  +3.46%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09070"
$ws.Range("E31").Value = "  +3.51%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.016"
$ws.Range("E32").Value = "  +0.38%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05204"
$ws.Range("E33").Value = "  +3.08%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7485"
$ws.Range("E34").Value = "  +4.11%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.110"
$ws.Range("E35").Value = "  -0.28%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.773"
$ws.Range("E36").Value = "  +4.32%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01837"
$ws.Range("E37").Value = "  -0.38%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.682"
$ws.Range("E38").Value = "  -0.44%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.145"
$ws.Range("E39").Value = "  -0.61%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9389"
$ws.Range("E40").Value = "  +0.85%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4423"
$ws.Range("E41").Value = "  +4.54%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "105.43"
$ws.Range("E42").Value = "  +2.02%  "
$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.000"
$ws.Range("E43").Value = "  +0.19%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.763"
$ws.Range("E44").Value = "  -0.14%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.607"
$ws.Range("E45").Value = "  +2.90%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1346"
$ws.Range("E46").Value = "  +5.59%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05870"
$ws.Range("E47").Value = "  +2.65%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.710"
$ws.Range("E48").Value = "  +4.69%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.424"
$ws.Range("E49").Value = "  +6.19%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.3938"
$ws.Range("E50").Value = "  +4.15%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "33.60"
$ws.Range("E51").Value = "  +2.48%  "
